# Apply updated cryptocurrency price/volume values to the active worksheet.
# This mirrors the scraper refresh commit "Updated cryptos list ... with GitHub Actions".
#
# For the Price column (D), values such as "0.9997" or "5.830" look numeric to
# Excel and would otherwise be auto-converted/normalized (losing trailing zeros,
# switching to scientific notation, etc). Force the cell to Text format first so
# the literal string is preserved exactly, then restore the default "Normal" style
# so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "29.616.93"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +2.53%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.859.93"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.86%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9997"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "245.30"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.09%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.6979"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("E7").Value = "  +0.08%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.07727"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.45%  "
$ws.Range("E9").Value = "  +1.14%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "23.66"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.90%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07754"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.09%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "5.162"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.45%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.857.70"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.80%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "92.34"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.45%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.6921"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +3.00%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "6.558"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +3.31%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "29.600.94"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +2.52%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008344"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.02%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "2.104.99"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.77%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "241.91"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.13%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "12.76"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.14%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.9998"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.07%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "7.616"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.94%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  +2.47%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "8.906"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.13%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "159.78"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.73%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "18.30"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.76%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.535"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.27%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.251"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.47%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.185"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("E33").Value = "  +0.33%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.7781"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.55%  "
$ws.Range("E35").Value = "  +4.89%  "
$ws.Range("E36").Value = "  +1.66%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.687"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.30%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.323.17"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +10.48%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01875"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +2.10%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.735"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +2.28%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.9577"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +2.78%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "106.34"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -1.79%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "5.830"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +12.26%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.9999"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.00000000125"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.79%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "9.775"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +3.45%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.004.82"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +1.55%  "
$ws.Range("E48").Value = "  +1.04%  "
$ws.Range("E49").Value = "  +3.52%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "64.42"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +4.42%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "6.986"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.86%  "
